$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "22.429.62"
$ws.Range("E2").Value = "  +0.17%  "
$ws.Range("D3").Value = "1.574.65"
$ws.Range("E3").Value = "  +0.50%  "
$ws.Range("D4").Value = "'1.003"
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").Value = "'1.002"
$ws.Range("E5").Value = "  +0.11%  "
$ws.Range("D6").Value = "'291.31"
$ws.Range("E6").Value = "  +0.00%  "
$ws.Range("D7").Value = "'0.3765"
$ws.Range("E7").Value = "  +2.58%  "
$ws.Range("D8").Value = "'49.90"
$ws.Range("E8").Value = "  +0.76%  "
$ws.Range("D9").Value = "'0.3423"
$ws.Range("D10").Value = "'1.164"
$ws.Range("E10").Value = "  -0.19%  "
$ws.Range("D11").Value = "'0.07679"
$ws.Range("E11").Value = "  +1.62%  "
$ws.Range("D12").Value = "'1.003"
$ws.Range("E12").Value = "  +0.15%  "
$ws.Range("E13").Value = "  +1.22%  "
$ws.Range("D14").Value = "'6.006"
$ws.Range("E14").Value = "  -0.53%  "
$ws.Range("D15").Value = "'6.934"
$ws.Range("E15").Value = "  +1.31%  "
$ws.Range("D16").Value = "1.574.85"
$ws.Range("E16").Value = "  -0.01%  "
$ws.Range("D17").Value = "'0.00001136"
$ws.Range("E17").Value = "  -0.50%  "
$ws.Range("D18").Value = "'90.39"
$ws.Range("E18").Value = "  +1.42%  "
$ws.Range("D19").Value = "'0.06762"
$ws.Range("E19").Value = "  +1.02%  "
$ws.Range("E20").Value = "  +0.12%  "
$ws.Range("E21").Value = "  +3.05%  "
$ws.Range("D22").Value = "'6.240"
$ws.Range("E22").Value = "  +0.01%  "
$ws.Range("D23").Value = "'12.07"
$ws.Range("E23").Value = "  +1.24%  "
$ws.Range("D24").Value = "'2.430"
$ws.Range("E24").Value = "  +1.35%  "
$ws.Range("D25").Value = "22.432.40"
$ws.Range("E25").Value = "  +0.16%  "
$ws.Range("D26").Value = "'2.744"
$ws.Range("E26").Value = "  -6.09%  "
$ws.Range("D27").Value = "'20.35"
$ws.Range("E27").Value = "  +2.69%  "
$ws.Range("D28").Value = "'146.02"
$ws.Range("E28").Value = "  -0.51%  "
$ws.Range("D29").Value = "'5.027"
$ws.Range("E29").Value = "  +1.66%  "
$ws.Range("D30").Value = "'126.32"
$ws.Range("E30").Value = "  +1.26%  "
$ws.Range("D31").Value = "1.749.66"
$ws.Range("E31").Value = "  -0.04%  "
$ws.Range("D32").Value = "'6.232"
$ws.Range("E32").Value = "  -0.29%  "
$ws.Range("D33").Value = "'1.012"
$ws.Range("E33").Value = "  +2.53%  "
$ws.Range("D34").Value = "'2.015"
$ws.Range("E34").Value = "  +1.90%  "
$ws.Range("D35").Value = "'10.03"
$ws.Range("E35").Value = "  -2.99%  "
$ws.Range("E36").Value = "  +1.88%  "
$ws.Range("D37").Value = "'0.02557"
$ws.Range("E37").Value = "  +1.31%  "
$ws.Range("E38").Value = "  +0.96%  "
$ws.Range("D39").Value = "'0.06580"
$ws.Range("E39").Value = "  +1.42%  "
$ws.Range("D40").Value = "'1.339"
$ws.Range("E40").Value = "  +7.87%  "
$ws.Range("D41").Value = "'5.463"
$ws.Range("E41").Value = "  -0.79%  "
$ws.Range("B42").Value = "TheSandbox"
$ws.Range("C42").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D42").Value = "'0.6469"
$ws.Range("E42").Value = "  +1.60%  "
$ws.Range("B43").Value = "Aptos"
$ws.Range("C43").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D43").Value = "'11.58"
$ws.Range("E43").Value = "  -1.49%  "
$ws.Range("D44").Value = "'14.23"
$ws.Range("E44").Value = "  -1.31%  "
$ws.Range("E45").Value = "  +0.18%  "
$ws.Range("D46").Value = "'0.6029"
$ws.Range("E46").Value = "  +0.15%  "
$ws.Range("D47").Value = "'3.802"
$ws.Range("E47").Value = "  +0.80%  "
$ws.Range("D48").Value = "'1.300"
$ws.Range("E48").Value = "  +10.08%  "
$ws.Range("D49").Value = "'2.090"
$ws.Range("E49").Value = "  -0.90%  "
$ws.Range("D50").Value = "'125.44"
$ws.Range("E50").Value = "  +3.26%  "
$ws.Range("D51").Value = "'0.07324"
$ws.Range("E51").Value = "  +0.82%  "
